$p = $ppt.ActivePresentation

# Remove the last slide (sldId 261 / slide6.xml) from the deck.
$s = $p.Slides.Item($p.Slides.Count)
$s.Delete()
